# daily auto push: 2026-02-11 10:05 UTC
# A new timestamp reading for 2026/02/11 (水) was appended to the log,
# landing right after the existing 2026/02/11 row (old row 797) and
# pushing every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 797 (and everything below it) down by one to make room for
# the new reading.
$ws.Rows.Item(797).Insert()

# Column A holds the date as plain text (e.g. "2026/02/11"). Assigning a
# date-shaped string straight to .Value would get auto-parsed into a date
# serial by Excel's type inference, which the source data does not use.
# Copy the date text from the row above instead, so it round-trips as the
# same literal string/type with no extra number formatting.
$ws.Cells.Item(796, 1).Copy()
$ws.Cells.Item(797, 1).PasteSpecial()

$ws.Cells.Item(797, 2).Value = "水"
$ws.Cells.Item(797, 3).Value = 17
$ws.Cells.Item(797, 4).Value = 201
